$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @{
    "D2" = "'243.22"
    "D3" = "'23.58"
    "D4" = "'5.284"
    "D5" = "'0.05785"
    "D6" = "'6.483"
    "D7" = "'3.338"
    "D8" = "'0.8092"
    "D9" = "'0.8766"
    "B10" = "One"
    "C10" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D10" = "'0.01035"
    "E10" = "9OneONEBestin24h"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D11" = "'0.1382"
    "E11" = "10WazirXWRX"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D12" = "'0.07276"
    "E12" = "11MandalaExchangeTokenMDX"
    "B13" = "LiechtensteinCryptoassetsExchange"
    "C13" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D13" = "'0.03091"
    "E13" = "12LiechtensteinCryptoassetsExchangeLCX"
    "B14" = "BitrueCoin"
    "C14" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D14" = "'0.03059"
    "E14" = "13BitrueCoinBTR"
    "B15" = "BitMartToken"
    "C15" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D15" = "'0.09316"
    "E15" = "14BitMartTokenBMX"
    "B16" = "MCDex"
    "C16" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "D16" = "'3.847"
    "E16" = "15MCDexMCB"
    "B17" = "BitForexToken"
    "C17" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D17" = "'0.001534"
    "E17" = "16BitForexTokenBF"
    "B18" = "CoinExToken"
    "C18" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "D18" = "'0.04705"
    "E18" = "17CoinExTokenCET"
    "D19" = "'0.006048"
    "D20" = "'0.001289"
    "D21" = "'0.004600"
    "D22" = "'0.00008701"
    "E22" = "21NitroExNTX"
    "D23" = "'3.578"
    "D24" = "'2.142"
    "D25" = "'0.3182"
    "D26" = "'0.1318"
    "D40" = "'0.03771"
    "D41" = "'0.006381"
    "B42" = "CEJI"
    "C42" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "D42" = "'0.004000"
    "E42" = "41CEJICEJI"
    "B43" = "BKEXToken"
    "C43" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D43" = "'0.1051"
    "E43" = "42BKEXTokenBKK"
    "D44" = "'0.007143"
    "E44" = "43LocalTradersLCTWorstin24h"
    "D45" = "'0.00005467"
    "D47" = "'0.5901"
    "D48" = "'0.001858"
}

foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).Value = $cellUpdates[$addr]
}
